$p = $ppt.ActivePresentation

# --- Slide 16: merge the "              " (14 spaces) run with the
#     following "－" run into a single run "              －" ---
$s16 = $p.Slides.Item(16)
$sh16 = $s16.Shapes.Item(2)
$tr16 = $sh16.TextFrame.TextRange
$dash = [char]0xFF0D
$merged = "              " + $dash
$c16 = $tr16.Characters(87, 15)
$c16.Text = $merged

# --- Slide 4: split the "最后 " run into three runs "最", "近", " " ---
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)
$tr4 = $sh4.TextFrame.TextRange
$c4 = $tr4.Characters(49, 1)
$c4.Text = "近"
